# This script rebuilds the "Sheet1" testcase table to match the updated
# schema: a new "Populasi" / "Generasi ACO" column layout (B..AA) and
# refreshed benchmark figures, per commit "feature: add column testcase".
#
# Net column changes (old -> new):
#   B Generasi          -> B Generasi GA
#   C Jumlah Semut      -> E Jumlah Semut   (moved right)
#   D Populasi          -> C Populasi       (moved left)
#   (new)               -> D Generasi ACO   (brand new column)
#   E..Z (Alpha..Image GA-ACO Percobaan 3) -> F..AA (shifted right by 1)
#
# Because several of the measured values (distances/runtimes) were also
# recomputed (not just shifted), the whole A1:AA10 block is written
# explicitly with its final values rather than relying on Excel's
# insert/move-column commands.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full 10-row x 27-column (A:AA) block of final values.
$data = New-Object 'object[,]' 10,27

$data[0,0] = 'Dataset'
$data[0,1] = 'Generasi GA'
$data[0,2] = 'Populasi'
$data[0,3] = 'Generasi ACO'
$data[0,4] = 'Jumlah Semut'
$data[0,5] = 'Alpha'
$data[0,6] = 'Beta'
$data[0,7] = 'Rho'
$data[0,8] = 'Pheromone Awal'
$data[0,9] = 'Jarak GA Percobaan 1'
$data[0,10] = 'Jarak GA Percobaan 2'
$data[0,11] = 'Jarak GA Percobaan 3'
$data[0,12] = 'Jarak GA-ACO Percobaan 1'
$data[0,13] = 'Jarak GA-ACO Percobaan 2'
$data[0,14] = 'Jarak GA-ACO Percobaan 3'
$data[0,15] = 'Runtime GA Percobaan 1'
$data[0,16] = 'Runtime GA Percobaan 2'
$data[0,17] = 'Runtime GA Percobaan 3'
$data[0,18] = 'Runtime GA-ACO Percobaan 1'
$data[0,19] = 'Runtime GA-ACO Percobaan 2'
$data[0,20] = 'Runtime GA-ACO Percobaan 3'
$data[0,21] = 'Image GA Percobaan 1'
$data[0,22] = 'Image GA Percobaan 2'
$data[0,23] = 'Image GA Percobaan 3'
$data[0,24] = 'Image GA-ACO Percobaan 1'
$data[0,25] = 'Image GA-ACO Percobaan 2'
$data[0,26] = 'Image GA-ACO Percobaan 3'
$data[1,0] = 't5.csv'
$data[1,1] = 10
$data[1,2] = 5
$data[1,3] = 10
$data[1,4] = 5
$data[1,5] = 1
$data[1,6] = 1
$data[1,7] = 0.5
$data[1,8] = 10
$data[1,9] = 102.1070184418289
$data[1,10] = 102.1070184418289
$data[1,11] = 102.1070184418289
$data[1,12] = 104.0130965649983
$data[1,13] = 104.0130965649983
$data[1,14] = 104.0130965649983
$data[1,15] = 0.0012
$data[1,16] = 0.0009
$data[1,17] = 0.0008
$data[1,18] = 0.0008
$data[1,19] = 0.0008
$data[1,20] = 0.0008
$data[1,21] = './imageResult/t5_1_GA_10.png'
$data[1,22] = './imageResult/t5_2_GA_10.png'
$data[1,23] = './imageResult/t5_3_GA_10.png'
$data[1,24] = './imageResult/t5_1_ACO_10.png'
$data[1,25] = './imageResult/t5_2_ACO_10.png'
$data[1,26] = './imageResult/t5_3_ACO_10.png'
$data[2,0] = 't5.csv'
$data[2,1] = 50
$data[2,2] = 5
$data[2,3] = 50
$data[2,4] = 5
$data[2,5] = 1
$data[2,6] = 1
$data[2,7] = 0.5
$data[2,8] = 10
$data[2,9] = 102.1070184418289
$data[2,10] = 102.1070184418289
$data[2,11] = 102.1070184418289
$data[2,12] = 104.0130965649983
$data[2,13] = 104.0130965649983
$data[2,14] = 104.0130965649983
$data[2,15] = 0.0009
$data[2,16] = 0.0009
$data[2,17] = 0.0009
$data[2,18] = 0.0014
$data[2,19] = 0.0009
$data[2,20] = 0.0009
$data[2,21] = './imageResult/t5_1_GA_50.png'
$data[2,22] = './imageResult/t5_2_GA_50.png'
$data[2,23] = './imageResult/t5_3_GA_50.png'
$data[2,24] = './imageResult/t5_1_ACO_50.png'
$data[2,25] = './imageResult/t5_2_ACO_50.png'
$data[2,26] = './imageResult/t5_3_ACO_50.png'
$data[3,0] = 't5.csv'
$data[3,1] = 100
$data[3,2] = 5
$data[3,3] = 100
$data[3,4] = 5
$data[3,5] = 1
$data[3,6] = 1
$data[3,7] = 0.5
$data[3,8] = 10
$data[3,9] = 102.1070184418289
$data[3,10] = 102.1070184418289
$data[3,11] = 102.1070184418289
$data[3,12] = 104.0130965649983
$data[3,13] = 104.0130965649983
$data[3,14] = 104.0130965649983
$data[3,15] = 0.0011
$data[3,16] = 0.0011
$data[3,17] = 0.0011
$data[3,18] = 0.001
$data[3,19] = 0.001
$data[3,20] = 0.001
$data[3,21] = './imageResult/t5_1_GA_100.png'
$data[3,22] = './imageResult/t5_2_GA_100.png'
$data[3,23] = './imageResult/t5_3_GA_100.png'
$data[3,24] = './imageResult/t5_1_ACO_100.png'
$data[3,25] = './imageResult/t5_2_ACO_100.png'
$data[3,26] = './imageResult/t5_3_ACO_100.png'
$data[4,0] = 'burma14.csv'
$data[4,1] = 10
$data[4,2] = 10
$data[4,3] = 10
$data[4,4] = 10
$data[4,5] = 1
$data[4,6] = 1
$data[4,7] = 0.5
$data[4,8] = 10
$data[4,9] = 52.29157121533663
$data[4,10] = 50.04909059225493
$data[4,11] = 51.47932488581549
$data[4,12] = 31.88252949105588
$data[4,13] = 31.88252949105588
$data[4,14] = 31.88252949105588
$data[4,15] = 0.0011
$data[4,16] = 0.001
$data[4,17] = 0.0011
$data[4,18] = 0.001
$data[4,19] = 0.0009
$data[4,20] = 0.0009
$data[4,21] = './imageResult/burma14_1_GA_10.png'
$data[4,22] = './imageResult/burma14_2_GA_10.png'
$data[4,23] = './imageResult/burma14_3_GA_10.png'
$data[4,24] = './imageResult/burma14_1_ACO_10.png'
$data[4,25] = './imageResult/burma14_2_ACO_10.png'
$data[4,26] = './imageResult/burma14_3_ACO_10.png'
$data[5,0] = 'burma14.csv'
$data[5,1] = 50
$data[5,2] = 10
$data[5,3] = 50
$data[5,4] = 10
$data[5,5] = 1
$data[5,6] = 1
$data[5,7] = 0.5
$data[5,8] = 10
$data[5,9] = 38.97701604838954
$data[5,10] = 44.13861991124494
$data[5,11] = 46.76581113716075
$data[5,12] = 31.45623383762054
$data[5,13] = 31.88252949105588
$data[5,14] = 31.88252949105588
$data[5,15] = 0.002
$data[5,16] = 0.0021
$data[5,17] = 0.0021
$data[5,18] = 0.0016
$data[5,19] = 0.0016
$data[5,20] = 0.0016
$data[5,21] = './imageResult/burma14_1_GA_50.png'
$data[5,22] = './imageResult/burma14_2_GA_50.png'
$data[5,23] = './imageResult/burma14_3_GA_50.png'
$data[5,24] = './imageResult/burma14_1_ACO_50.png'
$data[5,25] = './imageResult/burma14_2_ACO_50.png'
$data[5,26] = './imageResult/burma14_3_ACO_50.png'
$data[6,0] = 'burma14.csv'
$data[6,1] = 100
$data[6,2] = 10
$data[6,3] = 100
$data[6,4] = 10
$data[6,5] = 1
$data[6,6] = 1
$data[6,7] = 0.5
$data[6,8] = 10
$data[6,9] = 39.36589804668111
$data[6,10] = 40.79110455795923
$data[6,11] = 38.17739758246899
$data[6,12] = 31.88252949105588
$data[6,13] = 31.88252949105588
$data[6,14] = 31.22691510942754
$data[6,15] = 0.0034
$data[6,16] = 0.0039
$data[6,17] = 0.0034
$data[6,18] = 0.0024
$data[6,19] = 0.0023
$data[6,20] = 0.0023
$data[6,21] = './imageResult/burma14_1_GA_100.png'
$data[6,22] = './imageResult/burma14_2_GA_100.png'
$data[6,23] = './imageResult/burma14_3_GA_100.png'
$data[6,24] = './imageResult/burma14_1_ACO_100.png'
$data[6,25] = './imageResult/burma14_2_ACO_100.png'
$data[6,26] = './imageResult/burma14_3_ACO_100.png'
$data[7,0] = 'lin318.csv'
$data[7,1] = 10
$data[7,2] = 10
$data[7,3] = 10
$data[7,4] = 100
$data[7,5] = 1
$data[7,6] = 1
$data[7,7] = 0.5
$data[7,8] = 10
$data[7,9] = 557916.5971658916
$data[7,10] = 568259.1518398157
$data[7,11] = 564105.4600579566
$data[7,12] = 49215.61251916289
$data[7,13] = 49143.7729793856
$data[7,14] = 49215.61251916289
$data[7,15] = 0.0081
$data[7,16] = 0.008200000000000001
$data[7,17] = 0.008
$data[7,18] = 0.0668
$data[7,19] = 0.06560000000000001
$data[7,20] = 0.06560000000000001
$data[7,21] = './imageResult/lin318_1_GA_10.png'
$data[7,22] = './imageResult/lin318_2_GA_10.png'
$data[7,23] = './imageResult/lin318_3_GA_10.png'
$data[7,24] = './imageResult/lin318_1_ACO_10.png'
$data[7,25] = './imageResult/lin318_2_ACO_10.png'
$data[7,26] = './imageResult/lin318_3_ACO_10.png'
$data[8,0] = 'lin318.csv'
$data[8,1] = 50
$data[8,2] = 10
$data[8,3] = 50
$data[8,4] = 100
$data[8,5] = 1
$data[8,6] = 1
$data[8,7] = 0.5
$data[8,8] = 10
$data[8,9] = 516103.7735556596
$data[8,10] = 511781.1117244247
$data[8,11] = 537484.570807148
$data[8,12] = 49294.74163904427
$data[8,13] = 49215.61251916289
$data[8,14] = 48320.84193889733
$data[8,15] = 0.0337
$data[8,16] = 0.034
$data[8,17] = 0.0346
$data[8,18] = 0.3127
$data[8,19] = 0.3193
$data[8,20] = 0.3154
$data[8,21] = './imageResult/lin318_1_GA_50.png'
$data[8,22] = './imageResult/lin318_2_GA_50.png'
$data[8,23] = './imageResult/lin318_3_GA_50.png'
$data[8,24] = './imageResult/lin318_1_ACO_50.png'
$data[8,25] = './imageResult/lin318_2_ACO_50.png'
$data[8,26] = './imageResult/lin318_3_ACO_50.png'
$data[9,0] = 'lin318.csv'
$data[9,1] = 100
$data[9,2] = 10
$data[9,3] = 100
$data[9,4] = 100
$data[9,5] = 1
$data[9,6] = 1
$data[9,7] = 0.5
$data[9,8] = 10
$data[9,9] = 504397.0554574772
$data[9,10] = 516374.9896154903
$data[9,11] = 506839.1840932127
$data[9,12] = 49143.7729793856
$data[9,13] = 48272.74602537625
$data[9,14] = 49215.61251916289
$data[9,15] = 0.0672
$data[9,16] = 0.0665
$data[9,17] = 0.0667
$data[9,18] = 0.6264999999999999
$data[9,19] = 0.6267
$data[9,20] = 0.6378
$data[9,21] = './imageResult/lin318_1_GA_100.png'
$data[9,22] = './imageResult/lin318_2_GA_100.png'
$data[9,23] = './imageResult/lin318_3_GA_100.png'
$data[9,24] = './imageResult/lin318_1_ACO_100.png'
$data[9,25] = './imageResult/lin318_2_ACO_100.png'
$data[9,26] = './imageResult/lin318_3_ACO_100.png'

# Write the whole block in a single operation (A1:AA10).
$ws.Range("A1:AA10").Value = $data

# The newly introduced column AA has no formatting yet; copy the header
# cell formatting (bold font, borders, centered/top alignment) from the
# adjacent existing header cell Z1 so AA1 matches the rest of row 1.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("AA1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
